$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = "crop"
$ws.Range("B10").Value = "road_n_railway"
$ws.Range("B14").Value = "road_n_railway"
$ws.Range("B20").Value = "crop"
$ws.Range("B23").Value = "settlement"
$ws.Range("B24").Value = "agriculture"
$ws.Range("B80").Value = "river"
$ws.Range("B116").Value = "land_without_scrub"
$ws.Range("B152").Value = "crop"
$ws.Range("B156").Value = "crop"
$ws.Range("B168").Value = "land_without_scrub"
$ws.Range("B172").Value = "river"
$ws.Range("B174").Value = "forest"
$ws.Range("B176").Value = "crop"
$ws.Range("B182").Value = "river"
$ws.Range("B192").Value = "land_without_scrub"
$ws.Range("B194").Value = "land_without_scrub"
$ws.Range("B216").Value = "grassland"
$ws.Range("B217").Value = "crop"
$ws.Range("B242").Value = "grassland"
$ws.Range("B249").Value = "road_n_railway"
$ws.Range("B254").Value = "agriculture"
$ws.Range("B268").Value = "land_without_scrub"
$ws.Range("B269").Value = "land_without_scrub"
$ws.Range("B281").Value = "road_n_railway"
$ws.Range("B316").Value = "agriculture"
$ws.Range("B321").Value = "road_n_railway"
$ws.Range("B327").Value = "river"
$ws.Range("B356").Value = "grassland"
$ws.Range("B357").Value = "road_n_railway"
$ws.Range("B365").Value = "tank"
$ws.Range("B390").Value = "road_n_railway"
$ws.Range("B406").Value = "road_n_railway"
$ws.Range("B408").Value = "agriculture"
$ws.Range("B415").Value = "river"
$ws.Range("B423").Value = "forest"
$ws.Range("B432").Value = "agriculture"
$ws.Range("B451").Value = "land_without_scrub"
$ws.Range("B467").Value = "river"
$ws.Range("B471").Value = "road_n_railway"
$ws.Range("B477").Value = "road_n_railway"
$ws.Range("B481").Value = "road_n_railway"
$ws.Range("B502").Value = "agriculture"
$ws.Range("B506").Value = "river"
$ws.Range("B525").Value = "crop"
$ws.Range("B530").Value = "agriculture"
$ws.Range("B556").Value = "grassland"
$ws.Range("B571").Value = "land_without_scrub"
$ws.Range("B576").Value = "grassland"
$ws.Range("B583").Value = "forest"
$ws.Range("B591").Value = "crop"
$ws.Range("B597").Value = "road_n_railway"
$ws.Range("B605").Value = "agriculture"
$ws.Range("B608").Value = "river"
$ws.Range("B619").Value = "grassland"
$ws.Range("B624").Value = "agriculture"
$ws.Range("B626").Value = "land_without_scrub"
$ws.Range("B646").Value = "road_n_railway"
$ws.Range("B650").Value = "crop"
$ws.Range("B665").Value = "crop"
$ws.Range("B687").Value = "agriculture"
$ws.Range("B695").Value = "land_without_scrub"
$ws.Range("B709").Value = "road_n_railway"
$ws.Range("B711").Value = "crop"
$ws.Range("B729").Value = "agriculture"
$ws.Range("B745").Value = "road_n_railway"
$ws.Range("B757").Value = "agriculture"
$ws.Range("B775").Value = "crop"
$ws.Range("B784").Value = "road_n_railway"
$ws.Range("B800").Value = "road_n_railway"
$ws.Range("B810").Value = "crop"
$ws.Range("B842").Value = "land_without_scrub"
$ws.Range("B849").Value = "settlement"
$ws.Range("B854").Value = "crop"
$ws.Range("B872").Value = "land_without_scrub"
$ws.Range("B874").Value = "land_without_scrub"
$ws.Range("B880").Value = "land_without_scrub"
$ws.Range("B906").Value = "road_n_railway"
$ws.Range("B910").Value = "land_without_scrub"
$ws.Range("B948").Value = "grassland"
$ws.Range("B951").Value = "road_n_railway"
$ws.Range("B952").Value = "river"
$ws.Range("B1009").Value = "forest"
$ws.Range("B1014").Value = "river"
$ws.Range("B1023").Value = "crop"
$ws.Range("B1024").Value = "road_n_railway"
$ws.Range("B1026").Value = "crop"
$ws.Range("B1035").Value = "land_without_scrub"
$ws.Range("B1050").Value = "settlement"
$ws.Range("B1051").Value = "river"
$ws.Range("B1062").Value = "agriculture"
$ws.Range("B1063").Value = "crop"
$ws.Range("B1069").Value = "river"
$ws.Range("B1087").Value = "grassland"
$ws.Range("B1094").Value = "river"
$ws.Range("B1101").Value = "crop"
$ws.Range("B1104").Value = "agriculture"
$ws.Range("B1107").Value = "crop"
$ws.Range("B1111").Value = "river"
$ws.Range("B1117").Value = "agriculture"
$ws.Range("B1123").Value = "agriculture"
$ws.Range("B1135").Value = "crop"
$ws.Range("B1145").Value = "land_without_scrub"
$ws.Range("B1158").Value = "crop"
$ws.Range("B1159").Value = "crop"
$ws.Range("B1167").Value = "road_n_railway"
$ws.Range("B1182").Value = "crop"
$ws.Range("B1194").Value = "river"
$ws.Range("B1210").Value = "settlement"
$ws.Range("B1213").Value = "crop"
$ws.Range("B1222").Value = "agriculture"
$ws.Range("B1224").Value = "road_n_railway"
